# Auto-generated Word COM-interop script
# Replaces long paragraphs' single runs with multiple runs separated
# by manual line breaks (<w:br/>), per the target diff. Each Find.Execute
# targets a short, unique text window spanning one split point and
# re-inserts it with a Word manual-line-break code (^l) at that point.
# Replace:=1 (wdReplaceOne) so only the single targeted occurrence changes.
$d = $word.ActiveDocument

# prog_pt_1
$found = $d.Content.Find.Execute('rogenação, 5.4.3- Oxidação; 6-Indústria Cerâmica: 6.1 Ciment', $false, $false, $false, $false, $false, $true, 1, $false, 'rogenação, 5.4.3- Oxidação; 6-^lIndústria Cerâmica: 6.1 Ciment', 1)
if (-not $found) { Write-Output "FAILED: prog_pt_1" }

# prog_en_1
$found = $d.Content.Find.Execute('tion chain, 2.3- Chlorine-SodaSegment, 2.4- Synthesis Gas, M', $false, $false, $false, $false, $false, $true, 1, $false, 'tion chain, 2.3- Chlorine-Soda^lSegment, 2.4- Synthesis Gas, M', 1)
if (-not $found) { Write-Output "FAILED: prog_en_1" }

# prog_en_2
$found = $d.Content.Find.Execute('tion; 6- Ceramic Industry: 6.1Cement: 6.1.1- Overview, 6.1.2', $false, $false, $false, $false, $false, $true, 1, $false, 'tion; 6- Ceramic Industry: 6.1^lCement: 6.1.1- Overview, 6.1.2', 1)
if (-not $found) { Write-Output "FAILED: prog_en_2" }

# metodo_1
$found = $d.Content.Find.Execute('Método:Aulas expositivas, desenvolvim', $false, $false, $false, $false, $false, $true, 1, $false, 'Método:^lAulas expositivas, desenvolvim', 1)
if (-not $found) { Write-Output "FAILED: metodo_1" }

# metodo_2
$found = $d.Content.Find.Execute(', discussão de casos práticos.Critério:', $false, $false, $false, $false, $false, $true, 1, $false, ', discussão de casos práticos.^lCritério:', 1)
if (-not $found) { Write-Output "FAILED: metodo_2" }

# metodo_3
$found = $d.Content.Find.Execute('Critério:A nota (NOTA) será composta po', $false, $false, $false, $false, $false, $true, 1, $false, 'Critério:^lA nota (NOTA) será composta po', 1)
if (-not $found) { Write-Output "FAILED: metodo_3" }

# metodo_4
$found = $d.Content.Find.Execute(' a fórmula explicitada abaixo:NF = NOTA x % FREQ.', $false, $false, $false, $false, $false, $true, 1, $false, ' a fórmula explicitada abaixo:^lNF = NOTA x % FREQ.', 1)
if (-not $found) { Write-Output "FAILED: metodo_4" }

# metodo_5
$found = $d.Content.Find.Execute('NF = NOTA x % FREQ.Norma de Recuperação:', $false, $false, $false, $false, $false, $true, 1, $false, 'NF = NOTA x % FREQ.^lNorma de Recuperação:', 1)
if (-not $found) { Write-Output "FAILED: metodo_5" }

# metodo_6
$found = $d.Content.Find.Execute('Norma de Recuperação:Frequência mínima de 70% e not', $false, $false, $false, $false, $false, $true, 1, $false, 'Norma de Recuperação:^lFrequência mínima de 70% e not', 1)
if (-not $found) { Write-Output "FAILED: metodo_6" }

# criterio_1
$found = $d.Content.Find.Execute(' a fórmula explicitada abaixo:NF = NOTA x % FREQ.', $false, $false, $false, $false, $false, $true, 1, $false, ' a fórmula explicitada abaixo:^lNF = NOTA x % FREQ.', 1)
if (-not $found) { Write-Output "FAILED: criterio_1" }

# biblio_1
$found = $d.Content.Find.Execute('laborado pelo docente. Livros:Ullmann’s encyclopedia of indu', $false, $false, $false, $false, $false, $true, 1, $false, 'laborado pelo docente. Livros:^lUllmann’s encyclopedia of indu', 1)
if (-not $found) { Write-Output "FAILED: biblio_1" }

# biblio_2
$found = $d.Content.Find.Execute('m ; New York : WileyVCH, 2011.Encyclopedia of Chemical Proce', $false, $false, $false, $false, $false, $true, 1, $false, 'm ; New York : WileyVCH, 2011.^lEncyclopedia of Chemical Proce', 1)
if (-not $found) { Write-Output "FAILED: biblio_2" }

# biblio_3
$found = $d.Content.Find.Execute('York : Taylor & Francis, 2006.Manual Econômico da Indústria ', $false, $false, $false, $false, $false, $true, 1, $false, 'York : Taylor & Francis, 2006.^lManual Econômico da Indústria ', 1)
if (-not $found) { Write-Output "FAILED: biblio_3" }

# biblio_4
$found = $d.Content.Find.Execute('o; 8ed; Camaçari: CEPED, 2007.Shreve, R. Norris; BRINK JR., ', $false, $false, $false, $false, $false, $true, 1, $false, 'o; 8ed; Camaçari: CEPED, 2007.^lShreve, R. Norris; BRINK JR., ', 1)
if (-not $found) { Write-Output "FAILED: biblio_4" }

# biblio_5
$found = $d.Content.Find.Execute('o de Horácio Macedo; 4.ed. Riode Janeiro: Editora Guanabara ', $false, $false, $false, $false, $false, $true, 1, $false, 'o de Horácio Macedo; 4.ed. Rio^lde Janeiro: Editora Guanabara ', 1)
if (-not $found) { Write-Output "FAILED: biblio_5" }

# biblio_6
$found = $d.Content.Find.Execute('Koogan, 2008, c1997. Revistas:Química & Derivados, São Paulo', $false, $false, $false, $false, $false, $true, 1, $false, 'Koogan, 2008, c1997. Revistas:^lQuímica & Derivados, São Paulo', 1)
if (-not $found) { Write-Output "FAILED: biblio_6" }

# biblio_7
$found = $d.Content.Find.Execute('imica.com.br/category/revista/Petróleo & Energia, São Paulo,', $false, $false, $false, $false, $false, $true, 1, $false, 'imica.com.br/category/revista/^lPetróleo & Energia, São Paulo,', 1)
if (-not $found) { Write-Output "FAILED: biblio_7" }

# biblio_8
$found = $d.Content.Find.Execute('ry/revista-petroleo-e-energia/Revista FACTO, Publicação da A', $false, $false, $false, $false, $false, $true, 1, $false, 'ry/revista-petroleo-e-energia/^lRevista FACTO, Publicação da A', 1)
if (-not $found) { Write-Output "FAILED: biblio_8" }

# biblio_9
$found = $d.Content.Find.Execute('tp://www.abifina.org.br/facto/Revista Óleos & Gorduras, disp', $false, $false, $false, $false, $false, $true, 1, $false, 'tp://www.abifina.org.br/facto/^lRevista Óleos & Gorduras, disp', 1)
if (-not $found) { Write-Output "FAILED: biblio_9" }

Write-Output "DONE"
